$d = $word.ActiveDocument

# 1. Update the date/time line
$d.Content.Find.Execute("Ngày đặt : 2023/6/15. Thời gian : 18:1/16", $true, $false, $false, $false, $false, $true, 1, $false, "Ngày đặt : 2023/6/15. Thời gian : 21:54/20", 2)

# 2. Update table row 2 (Asus 998 -> ốp lưng 998, price change)
$t = $d.Tables.Item(1)
$t.Cell(2, 2).Range.Text = "ốp lưng 998"
$t.Cell(2, 3).Range.Text = "685,621 vnđ vnd"

# 3. Update table row 3 (Xiaomi 999 -> Samsung 999, price change)
$t.Cell(3, 2).Range.Text = "Samsung 999"
$t.Cell(3, 3).Range.Text = "1,375,101 vnđ vnd"

# 4. Add a new row for "Apple (iPad) 1000"
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Split(1, 4)
$newRow.Cells.Item(2).Split(1, 2)
$newRow.Cells.Item(3).Split(1, 2)
$newRow.Cells.Item(1).Range.Text = "3"
$newRow.Cells.Item(2).Range.Text = "Apple (iPad) 1000"
$newRow.Cells.Item(3).Range.Text = "10,269,329 vnđ vnd"
$newRow.Cells.Item(4).Range.Text = "1"

# 5. Update the total amount line
$d.Content.Find.Execute("Tổng tiền : 17,118,932 vnđ vnd", $true, $false, $false, $false, $false, $true, 1, $false, "Tổng tiền : 12,330,051 vnđ vnd", 2)
